# Atualização de bases das ligas, do dia: 19-06-2024 às 21:51
#
# The source data rows got re-synced with the upstream feed, which caused a
# handful of fixture rows to swap/rotate their match data (everything except
# the leading row-index column A) with neighboring rows.
#
# Row pairs/groups that exchange their B:AD content:
#   95  <-> 96            (simple swap)
#   129 -> 130 -> 131 -> 129   (3-way rotation)
#   224 <-> 225           (simple swap)
#   256 <-> 257           (simple swap)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple two-row swaps -------------------------------------------------
function Swap-Rows($ws, $r1, $r2) {
    $rng1 = $ws.Range("B$r1`:AD$r1")
    $rng2 = $ws.Range("B$r2`:AD$r2")
    $v1 = $rng1.Value2
    $v2 = $rng2.Value2
    $rng1.Value = $v2
    $rng2.Value = $v1
}

Swap-Rows $ws 95 96
Swap-Rows $ws 224 225
Swap-Rows $ws 256 257

# --- Three-way rotation: 129 -> 130 -> 131 -> 129 -------------------------
# New129 = Old131 ; New130 = Old129 ; New131 = Old130
$rng129 = $ws.Range("B129:AD129")
$rng130 = $ws.Range("B130:AD130")
$rng131 = $ws.Range("B131:AD131")

$v129 = $rng129.Value2
$v130 = $rng130.Value2
$v131 = $rng131.Value2

$rng129.Value = $v131
$rng130.Value = $v129
$rng131.Value = $v130
